$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 3.034748368925986

$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 6.741336633845642

$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 9.226618575922256
$ws.Range("D4").Value = 16.98373111632243
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 35.87465599807698

$ws.Range("B5").Value = 1.505614041169197
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 3.811642989160245

$ws.Range("B6").Value = 3.182878228561681
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.1529057820181812
$ws.Range("E6").Value = 6.48142807727062
$ws.Range("G6").Value = 11.47044854674929

$ws.Range("B7").Value = 3.182878228561681
$ws.Range("C7").Value = 1.65323645889881
$ws.Range("D7").Value = 0.7127328510149897
$ws.Range("E7").Value = 6.48142807727062
$ws.Range("G7").Value = 12.0302756157461

$ws.Range("B8").Value = 0.7287194209349384
$ws.Range("C8").Value = 1.65323645889881
$ws.Range("D8").Value = 0.7127328510149897
$ws.Range("E8").Value = 0.4998867070740569
$ws.Range("G8").Value = 3.594575437922795

$ws.Range("B9").Value = 3.182878228561681
$ws.Range("C9").Value = 1.65323645889881
$ws.Range("D9").Value = 0.7127328510149897
$ws.Range("E9").Value = 0.4998867070740569
$ws.Range("G9").Value = 6.048734245549538

$ws.Range("B10").Value = 3.182878228561681
$ws.Range("C10").Value = 1.65323645889881
$ws.Range("D10").Value = 0.7127328510149897
$ws.Range("E10").Value = 0.4998867070740569
$ws.Range("G10").Value = 6.048734245549538

$ws.Range("B11").Value = 0.0003714022599530242
$ws.Range("C11").Value = 0.3375848360084654
$ws.Range("D11").Value = 0.1529057820181812
$ws.Range("E11").Value = 0.4998867070740569
$ws.Range("G11").Value = 0.9907487273606566

$ws.Range("B12").Value = 3.182878228561681
$ws.Range("C12").Value = 1.65323645889881
$ws.Range("D12").Value = 16.98373111632243
$ws.Range("E12").Value = 6.48142807727062
$ws.Range("G12").Value = 28.30127388105354

$ws.Range("B13").Value = 0.7287194209349384
$ws.Range("C13").Value = 1.65323645889881
$ws.Range("D13").Value = 3.082599426703578
$ws.Range("E13").Value = 0.4998867070740569
$ws.Range("G13").Value = 5.964442013611383

$ws.Range("B14").Value = 0.7287194209349384
$ws.Range("C14").Value = 1.65323645889881
$ws.Range("D14").Value = 3.082599426703578
$ws.Range("E14").Value = 0.4998867070740569
$ws.Range("G14").Value = 5.964442013611383

$ws.Range("B15").Value = 3.182878228561681
$ws.Range("C15").Value = 1.65323645889881
$ws.Range("D15").Value = 16.98373111632243
$ws.Range("E15").Value = 6.48142807727062
$ws.Range("G15").Value = 28.30127388105354

$ws.Range("B16").Value = 0.1554434735375247
$ws.Range("C16").Value = 0.3375848360084654
$ws.Range("D16").Value = 3.082599426703578
$ws.Range("E16").Value = 0.4998867070740569
$ws.Range("G16").Value = 4.075514443323626

$ws.Range("B17").Value = 0.7287194209349384
$ws.Range("C17").Value = 0.05231270169004087
$ws.Range("D17").Value = 0.7127328510149897
$ws.Range("E17").Value = 0.4998867070740569
$ws.Range("G17").Value = 1.993651680714026

$ws.Range("B18").Value = 0.1554434735375247
$ws.Range("C18").Value = 0.3375848360084654
$ws.Range("D18").Value = 0.7127328510149897
$ws.Range("E18").Value = 0.4998867070740569
$ws.Range("G18").Value = 1.705647867635037

